$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "list of APIs" sheet: add a new row for the /signin route
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("list of APIs")
$lo  = $ws1.ListObjects.Item(1)

# Copy the formatting of the existing data row (row 2) down into row 3
# so the new row picks up the same fonts / borders used by the table.
$ws1.Range("A2:C2").Copy()
$ws1.Range("A3:C3").PasteSpecial(-4122)

$ws1.Range("A3").Value = "/signin"
$ws1.Range("B3").Value = "email, password"
$ws1.Range("C3").Value = "signin"

# ---------------------------------------------------------------------
# 2) Add the new "signin" worksheet, right after "signup"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("signup")
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "signin"

# ---------------------------------------------------------------------
# 3) Fill in the distinctive "Response" text for each table first
#    (this mirrors the order the strings were actually authored in)
# ---------------------------------------------------------------------
$ws3.Range("C4").Value = "`$user(data)"
$ws3.Range("C3").Value = "invalid password"
$ws2.Range("C4").Value = "Registration Succeeded"
$ws3.Range("C5").Value = "email not registered"

# Generic / repeated labels
$ws2.Range("B2").Value = "string"
$ws2.Range("A1").Value = "Output"
$ws2.Range("B1").Value = "Type"
$ws2.Range("C1").Value = "Response"

$ws3.Range("B4").Value = "array"

# ---------------------------------------------------------------------
# 4) Finish populating both response tables
# ---------------------------------------------------------------------
# signup
$ws2.Range("C2").Value = "invalid email"
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "string"
$ws2.Range("C3").Value = "email already registered"
$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "string"

# signin
$ws3.Range("A1").Value = "Output"
$ws3.Range("B1").Value = "Type"
$ws3.Range("C1").Value = "Response"
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "string"
$ws3.Range("C2").Value = "invalid email"
$ws3.Range("A3").Value = 2
$ws3.Range("A4").Value = 3
$ws3.Range("A5").Value = 4
$ws3.Range("B5").Value = "string"

# ---------------------------------------------------------------------
# 5) Header formatting (centered) for both response tables
# ---------------------------------------------------------------------
$ws2.Range("A1:C1").HorizontalAlignment = -4108
$ws2.Range("A2:A4").HorizontalAlignment = -4108
$ws2.Range("B2:B4").HorizontalAlignment = -4108

$ws3.Range("A1:C1").HorizontalAlignment = -4108
$ws3.Range("A2:A5").HorizontalAlignment = -4108
$ws3.Range("B2:B5").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 6) Hyperlink the new /signin row's "data" cell to the new sheet
# ---------------------------------------------------------------------
$ws1.Range("D3").Value = "data"
$ws1.Range("D3").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("D3"), "", "signin!A1", "", "data")

# Resize the table to include the new row
$lo.Resize($ws1.Range("A1:D3"))

# ---------------------------------------------------------------------
# 7) Selection tidy-up to mirror the final authored state
# ---------------------------------------------------------------------
$ws1.Range("D2").Select()
$ws3.Range("A1").Select()

Write-Host "done"
